# Update the Price (column D) and Volume(1h) (column E) values for the
# cryptos list, as produced by the scheduled GitHub Actions refresh on
# Sat May 13 19:48:21 UTC 2023. Every cell in D/E is stored as plain text
# in the workbook (see "Price"/"Volume(1h)" columns), including D values
# that look like plain numbers (e.g. "319.61") or have significant
# trailing zeros (e.g. "1.030"). Assigning such a string straight to
# .Value lets Excel's own type-inference silently reinterpret it as a
# number (dropping formatting / trailing zeros), so for every D-column
# price we force a text number-format first, assign the literal string,
# and then drop back to the default "Normal" style so no stray
# formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.522.56'
$ws.Range("E2").Value = '  +4.88%  '

$ws.Range("D3").Value = '1.841.65'
$ws.Range("E3").Value = '  +3.64%  '

$ws.Range("E4").Value = '  +2.63%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.61'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.98%  '

$ws.Range("E6").Value = '  +2.49%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4379'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.80%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3742'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.34%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07381'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8758'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.52%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.45'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.26%  '

$ws.Range("D12").Value = '1.855.40'
$ws.Range("E12").Value = '  +7.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.489'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.28%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.672'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07162'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.73'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.90%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.034'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.55%  '

$ws.Range("E18").Value = '  +3.64%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.027'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.41'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.34%  '

$ws.Range("D21").Value = '27.539.65'
$ws.Range("E21").Value = '  +4.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.257'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.79%  '

$ws.Range("E23").Value = '  +1.80%  '

$ws.Range("D24").Value = '2.060.49'
$ws.Range("E24").Value = '  +6.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.52'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.925'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.73'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.255'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.938'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.15'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09095'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.206'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7656'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.499'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.879'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.030'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.88%  '

$ws.Range("E37").Value = '  +3.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01975'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05259'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5167'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.789'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1668'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.657'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.511'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '108.88'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.51'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.031'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.706'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.11%  '

$ws.Range("E49").Value = '  +4.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06350'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.892'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.74%  '
